$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2:E23 values from 50 to 70
$ws.Range("E2:E23").Value = 70

# Update the selection to E30 (single cell)
$ws.Range("E30").Select()
